$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New table of attributes (eid_instance_id, eid_instance_attribute, datatype, profile_id, display_name)
# replacing the previous set of "unneeded instance attributes".
$data = @(
    @(4, "SEGMENT1", "mdex:string", 1, "Project Number"),
    @(4, "NAME", "mdex:string", 1, "Project Name"),
    @(4, "ORDER_NUMBER", "mdex:string", 1, "Order Number"),
    @(4, "ORDER_BOOK_DATE", "mdex:dateTime", 6, "Order Bridged Date"),
    @(4, "WWAPC", "mdex:string", 1, "WWAPC"),
    @(4, "SO_CURRENCY_CODE", "mdex:string", 1, "SO Currency"),
    @(4, "REVENUE_VALUE", "mdex:double", 6, "Revenue Value"),
    @(4, "RELEASED_AMOUNT", "mdex:double", 6, "Shipped Revenue"),
    @(4, "BACKLOG", "mdex:double", 6, "Backlog"),
    @(4, "FORECAST_EQP_COST", "mdex:double", 6, "Forecast EQP Cost"),
    @(4, "SHIPPED_FCST_EQP_COST", "mdex:double", 6, "Shipped Fcst EQP Cost"),
    @(4, "PERIOD_NAME", "mdex:string", 1, "Period"),
    @(4, "ORG_ID", "mdex:string", 1, "OU Name"),
    @(4, "CARRYING_OUT_ORGANIZTION_ID", "mdex:string", 1, "Organization")
)

$rowCount = $data.Length
$lastRow = 1 + $rowCount

for ($i = 0; $i -lt $rowCount; $i++) {
    $r = 2 + $i
    $row = $data[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
}

# Remove the now-unused trailing rows (previously rows 16-21).
$ws.Rows("16:21").Delete()

# Reset the view: clear the scrolled top-left cell and move the selection.
$ws.Range("C22").Select()
